$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of column P (the 2019 column) into the new column Q
# (the 2020 column) for every data row, then fill in the 2020 values.
# Using Copy + PasteSpecial(Formats) reproduces the same effective
# rendering (font, borders, alignment, number format) as column P.
$ws.Range("P3:P34").Copy()
$ws.Range("Q3:Q34").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Header (row 4) - new year column
$ws.Range("Q4").Value = 2020

# Data rows: value 2020 figures per indicator / region row.
# "-" (dash) denotes "no data" matching the rest of the table.
$ws.Range("Q5").Value = 51
$ws.Range("Q6").Value = 29
$ws.Range("Q7").Value = 22
$ws.Range("Q8").Value = 5
$ws.Range("Q9").Value = 3
$ws.Range("Q10").Value = 2
$ws.Range("Q11").Value = 15
$ws.Range("Q12").Value = 9
$ws.Range("Q13").Value = 5
$ws.Range("Q14").Value = "-"
$ws.Range("Q15").Value = "-"
$ws.Range("Q16").Value = "-"
$ws.Range("Q17").Value = "-"
$ws.Range("Q18").Value = "-"
$ws.Range("Q19").Value = "-"
$ws.Range("Q20").Value = 7
$ws.Range("Q21").Value = 7
$ws.Range("Q22").Value = "-"
$ws.Range("Q23").Value = "-"
$ws.Range("Q24").Value = "-"
$ws.Range("Q25").Value = "-"
$ws.Range("Q26").Value = 24
$ws.Range("Q27").Value = 10
$ws.Range("Q28").Value = 14
$ws.Range("Q29").Value = "-"
$ws.Range("Q30").Value = "-"
$ws.Range("Q31").Value = "-"
$ws.Range("Q32").Value = "-"
$ws.Range("Q33").Value = "-"
$ws.Range("Q34").Value = "-"

# Update the active selection to match the authored workbook state.
$ws.Range("K18").Select()
